# Updated symbol list on Sat Dec 24 05:31:09 UTC 2022 with GitHub Actions
# Refreshes crypto price (column D) / volume label (column E) data on Sheet1.
#
# Price cells hold numeric-looking text (e.g. "245.24") that must stay plain
# text, matching the source data file. A leading apostrophe forces Excel to
# keep the entry as text instead of auto-converting it to a number; resetting
# the cell Style back to "Normal" afterwards clears the quote-prefix
# formatting flag so no stray number-format/style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.Value = "'" + $val
    $rng.Style = "Normal"
}

Set-TextValue "D2"  "245.24"
Set-TextValue "D3"  "22.03"
Set-TextValue "D4"  "5.337"
Set-TextValue "D5"  "0.05972"
Set-TextValue "D6"  "3.400"
Set-TextValue "D7"  "6.388"
Set-TextValue "D8"  "0.8130"
Set-TextValue "D9"  "0.9665"
Set-TextValue "D10" "0.1428"
Set-TextValue "D12" "0.07407"
Set-TextValue "D13" "0.03056"
Set-TextValue "D14" "0.09404"
Set-TextValue "D15" "3.998"
Set-TextValue "D16" "0.001591"
Set-TextValue "D17" "0.04810"
Set-TextValue "D18" "0.0005915"
$ws.Range("E18").Value = "17OneONE"
Set-TextValue "D19" "0.006216"
Set-TextValue "D20" "0.004138"
Set-TextValue "D21" "0.0009868"
Set-TextValue "D22" "0.00009709"
Set-TextValue "D23" "3.745"
Set-TextValue "D24" "2.128"
Set-TextValue "D26" "0.1332"
Set-TextValue "D40" "0.03913"
Set-TextValue "D41" "0.006486"
Set-TextValue "D43" "0.003003"
Set-TextValue "D44" "0.005372"
Set-TextValue "D45" "0.00005314"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "D47" "0.8507"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
Set-TextValue "D48" "0.04028"
Set-TextValue "D49" "0.00002102"
